$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "Niccoló Orsi"
$ws.Range("B26").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C26").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("D26").Value = "Nicholas Marzadro | SBARX"
$ws.Range("E26").Value = "Marco Sartorelli | Modium"
$ws.Range("F26").Value = "Mattia Tezzele | U.SGUARNA"
